# Trade #28 closed at 2026-02-16 22:54:53 - base_strategy UP +0.000%
# Append the new trade row (row 29) to both the "All Trades" and the
# "base_strategy" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A29").Value = 28
    $ws.Range("B29").Value = "'2026-02-16"
    $ws.Range("C29").Value = "22:54:53"
    $ws.Range("D29").Value = "base_strategy"
    $ws.Range("E29").Value = "UP"
    $ws.Range("F29").Value = 49.999998
    $ws.Range("G29").Value = "'"
    $ws.Range("H29").Value = "OPEN"
    $ws.Range("I29").Value = 0
    $ws.Range("J29").Value = 0
    $ws.Range("K29").Value = 100
    $ws.Range("L29").Value = 0
    $ws.Range("M29").Value = 0
    $ws.Range("N29").Value = 0.6
    $ws.Range("O29").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P29").Value = "'"
    $ws.Range("Q29").Value = 0
}
